$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values are stored as text in the sheet (t="inlineStr").
# Force text formatting before assignment so Excel does not coerce these
# numeric-looking strings into Number cells (which would also introduce
# floating point / scientific-notation artifacts for values like 0.00000000751).
$priceUpdates = @{
    "D2" = "234.89"
    "D3" = "22.29"
    "D4" = "5.394"
    "D5" = "0.05639"
    "D6" = "6.486"
    "D7" = "3.374"
    "D8" = "1.068"
    "D9" = "0.7877"
    "D10" = "0.1401"
    "D11" = "0.07336"
    "D12" = "0.03209"
    "D13" = "0.02947"
    "D14" = "0.09256"
    "D15" = "0.001661"
    "D16" = "3.260"
    "D17" = "0.04758"
    "D18" = "0.0005747"
    "D19" = "0.006218"
    "D20" = "0.005095"
    "D21" = "0.001052"
    "D23" = "3.859"
    "D26" = "0.1054"
    "D27" = "0.0004996"
    "D40" = "0.04094"
    "D41" = "0.006966"
    "D42" = "0.1036"
    "D43" = "0.003248"
    "D44" = "0.009939"
    "D45" = "0.00005424"
    "D46" = "0.00000000751"
    "D47" = "0.6761"
    "D48" = "0.03879"
    "D49" = "0.00002103"
    "D50" = "0.01011"
}

foreach ($cellRef in $priceUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$cellRef]
    $cell.Style = "Normal"
}

# Column E ("Volume(1h)") label text updates.
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("E48").Value = "47BOLOBOLO"
